$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for column F (dSF) - row => value
$updates = @{
    2  = 1
    3  = -2
    5  = 6
    6  = -2
    7  = 1
    8  = 4
    9  = 1
    10 = 2
    11 = 7
    12 = 1
    15 = -2
    16 = 1
    17 = 1
    18 = 3
    19 = 2
    21 = -5
    22 = 2
    23 = 3
    24 = 3
    25 = 4
    26 = -2
    28 = 2
    29 = 1
    30 = 1
    31 = 3
    32 = 2
    33 = 1
    34 = -6
    35 = 2
    36 = 1
    37 = 2
    38 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
